$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Objetivos:): B/C value replaced by the docente string
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# Row 13 (was blank in A, now "Programa resumido:" with B/C "Semestral")
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 14 becomes "Short syllabus:" only (A), clear B/C
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# Row 15 becomes "Programa:" with B/C = "01/01/2021" (reuse shared string idx 14 via
# formula+paste-special trick so Excel doesn't auto-convert the text to a date serial,
# and so no new number-format style gets created).
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Formula = "=""01/01/2021"""
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").Formula = "=""01/01/2021"""
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Row 16 becomes "Syllabus:" only (A), clear B/C
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

# Row 17 becomes "Avaliação:" only (A)
$ws.Range("A17").Value = "Avaliação:"

# Row 18 becomes "Método:" with B/C = docente string
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"

# Row 19 becomes "Critério:" with B/C = "Aulas expositivas..."
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"

# Row 20 becomes "Norma de recuperação:" with B/C = "Média Aritmética..."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

# Row 21 becomes "Bibliografia:" with B/C = "NF = (MF + PR)/2, ..."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# Row 22 (old Bibliografia:/BAZZO... row) is removed entirely - delete the whole row
# so the sheet's used range shrinks back to A1:C21.
$ws.Range("A22:C22").EntireRow.Delete()
